# RBA v2.5 - Atualizacao da Tela
#
# Replaces the placeholder "tre"-family tokens (in various casings) with
# "qwer"-family tokens, both in the body text and in the primary header.
#
# wdReplaceOne (1) restricts Find.Execute to a single match so that each
# occurrence in a run of duplicate placeholders can receive its own
# (slightly different) replacement text.

$d = $word.ActiveDocument

# --- Body: bold "TERE" -> "QWER" ("A TERE," salutation line) ---------------
$bodyRng = $d.Content
$bodyRng.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# --- Primary header: sequence of "tre" placeholders -------------------------
$hdr = $d.Sections(1).Headers(1)

# Each entry: the (case-sensitive) text to search for, and its replacement,
# applied strictly in document order (one match consumed per entry).
$edits = @(
    @{ Find = "TRE";  Replace = "QWER" },
    @{ Find = "TERE"; Replace = "QWER" },
    @{ Find = "Tre";  Replace = "Qwer" },
    @{ Find = "Tre";  Replace = "Qwer" },
    @{ Find = "Tre";  Replace = "Qewr" },
    @{ Find = "Tre";  Replace = "Qewr" },
    @{ Find = "Tre";  Replace = "Qwer" },
    @{ Find = "tre";  Replace = "qwer" },
    @{ Find = "tre";  Replace = "qwer" },
    @{ Find = "tre";  Replace = "qwer" }
)

$rng = $hdr.Range
$rng.Start = 0
$rng.End = $hdr.Range.End

foreach ($edit in $edits) {
    $rng.Find.Execute($edit.Find, $true, $false, $false, $false, $false, $true, 1, $false, $edit.Replace, 1) | Out-Null
    $rng.Collapse(0)
    $rng.End = $hdr.Range.End
}
